$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("July 2018")
$ws.Activate()

# ---------------------------------------------------------------------------
# 1) Row 4 (existing "EVHP Concept 5 TEE Track Wide Bore" request): the print
#    finally completed, so fill in the "Date Completed" cell (B4).
#    "04-08-2018" is ambiguous (day=04, month=08 are both <=12) so Excel's
#    automatic data-type detection would otherwise turn it into a date
#    serial number instead of keeping it as literal text. To avoid that we
#    stage the literal text in a scratch cell that has been forced to Text
#    format, copy it, and paste *values only* into B4 - this keeps B4's own
#    formatting (and therefore its shared style) completely untouched.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "04-08-2018"
$scratch.Copy()
$ws.Range("B4").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
$scratch.ClearContents()
$scratch.Style = "Normal"
# B4 previously had its own directly-assigned (non-text) style; re-touching a
# harmless, already-true font property makes the cell resolve back to the
# same plain centered style used by its neighbours in the row.
$ws.Range("B4").Font.Name = "Calibri"

# ---------------------------------------------------------------------------
# 2) Row 5: brand-new request for the experimental "Spine in a Box 3lvl
#    Spine" part (date requested 16-07-2018, material Bridge, comments NA).
#    Column B (Date Completed) stays blank since the print hasn't finished.
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "16-07-2018"
$ws.Range("C5").Value = "Spine in a Box 3lvl Spine"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "Bridge"
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 20
$ws.Range("H5").Value = 0.2
$ws.Range("I5").Value = "NA"

# ---------------------------------------------------------------------------
# 3) Leave the cursor where the author last clicked after logging the row.
# ---------------------------------------------------------------------------
$ws.Range("F6").Select()
